$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("B2").Value = 0.5532819558862609
$ws.Range("D2").Value = "20210130"
$ws.Range("I2").Value = "50%"
$ws.Range("K2").Value = 2

# Update row 3 values
$ws.Range("B3").Value = 0.5500495540138751
$ws.Range("D3").Value = "20210130"
$ws.Range("I3").Value = "40%"
$ws.Range("K3").Value = 1

# Delete rows 4 and 5 (old rows no longer present)
$ws.Rows("4:5").Delete()
